$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Final review: update the Final-exam score for student 1 (H8).
# Dependent formulas (H7 weighted score, J7 total) recalc automatically.
$ws.Range("H8").Value = 0.71

# Leave the selection on H9, matching where the reviewer ended up.
$ws.Range("H9").Select()
